$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,10
$arr[0,0] = 0.1526000130636476
$arr[0,1] = 0.09915332967504753
$arr[0,2] = 0.140549707625361
$arr[0,3] = 2.491299424789091
$arr[0,4] = 1.779943708681913
$arr[0,5] = 1.57136609747451
$arr[0,6] = 1.775815328987598
$arr[0,7] = 0.21281199231926
$arr[0,8] = 2.62427832058529
$arr[0,9] = 0.2148883173652152
$arr[1,0] = 0.1490642734085696
$arr[1,1] = 0.09595602584956708
$arr[1,2] = 0.1394040164962149
$arr[1,3] = 2.510866881906729
$arr[1,4] = 1.798341471978617
$arr[1,5] = 1.587801569165379
$arr[1,6] = 1.790624223074772
$arr[1,7] = 0.2131302842509228
$arr[1,8] = 2.430284566814805
$arr[1,9] = 0.2145085662648327
$arr[2,0] = 0.1469291304647697
$arr[2,1] = 0.09400193139301649
$arr[2,2] = 0.1387456275885093
$arr[2,3] = 2.524486748038889
$arr[2,4] = 1.811006179643044
$arr[2,5] = 1.598795043588282
$arr[2,6] = 1.800882374654634
$arr[2,7] = 0.2134177468574698
$arr[2,8] = 2.311489827715889
$arr[2,9] = 0.2143525854560409
$arr[3,0] = 0.1460681149015812
$arr[3,1] = 0.0932079595699804
$arr[3,2] = 0.138488694624801
$arr[3,3] = 2.530440063428514
$arr[3,4] = 1.816510436888208
$arr[3,5] = 1.603501575869885
$arr[3,6] = 1.805355265272318
$arr[3,7] = 0.2135580431114086
$arr[3,8] = 2.263162216242449
$arr[3,9] = 0.2143084624953318
$arr[4,0] = 0.1459256938582598
$arr[4,1] = 0.09307626408217118
$arr[4,2] = 0.1384467185424967
$arr[4,3] = 2.531452936976237
$arr[4,4] = 1.817445124512986
$arr[4,5] = 1.604296771760801
$arr[4,6] = 1.806115645265741
$arr[4,7] = 0.2135827379147948
$arr[4,8] = 2.255142484092573
$arr[4,9] = 0.214302310989666
$arr[5,0] = 0.1469174816855627
$arr[5,1] = 0.0939912140673016
$arr[5,2] = 0.1387421164367382
$arr[5,3] = 2.524565405130346
$arr[5,4] = 1.811079023114516
$arr[5,5] = 1.598857600326198
$arr[5,6] = 1.800941513453552
$arr[5,7] = 0.2134195451813206
$arr[5,8] = 2.31083772884017
$arr[5,9] = 0.2143519116428294
$arr[6,0] = 0.1513734947317857
$arr[6,1] = 0.09804905202390302
$arr[6,2] = 0.1401453378571418
$arr[6,3] = 2.497712742518672
$arr[6,4] = 1.786002756131083
$arr[6,5] = 1.57684569980448
$arr[6,6] = 1.780679280546728
$arr[6,7] = 0.2129026447467481
$arr[6,8] = 2.557324467618116
$arr[6,9] = 0.2147413767360717
$arr[7,0] = 0.1603935584798535
$arr[7,1] = 0.1060762659699606
$arr[7,2] = 0.1432534999719834
$arr[7,3] = 2.457819697630455
$arr[7,4] = 1.747723062074698
$arr[7,5] = 1.540847639974828
$arr[7,6] = 1.750213178899124
$arr[7,7] = 0.2126189356680683
$arr[7,8] = 3.043140026529898
$arr[7,9] = 0.2161164627363092
$arr[8,0] = 0.167190109950738
$arr[8,1] = 0.112014276170882
$arr[8,2] = 0.1457532078503583
$arr[8,3] = 2.436331353252825
$arr[8,4] = 1.72629284749938
$arr[8,5] = 1.518783205885327
$arr[8,6] = 1.733509867280105
$arr[8,7] = 0.2128554019466478
$arr[8,8] = 3.401513353167388
$arr[8,9] = 0.2174983237112329
$arr[9,0] = 0.1703184285537844
$arr[9,1] = 0.1147240147617055
$arr[9,2] = 0.1469371111928446
$arr[9,3] = 2.428262661377076
$arr[9,4] = 1.718008526879316
$arr[9,5] = 1.509700412003838
$arr[9,6] = 1.727151090372537
$arr[9,7] = 0.2130595970184501
$arr[9,8] = 3.564852600194797
$arr[9,9] = 0.2182074535504412
$arr[10,0] = 0.1715082447064873
$arr[10,1] = 0.1157513007893698
$arr[10,2] = 0.1473921259732478
$arr[10,3] = 2.425453338651963
$arr[10,4] = 1.715082972680037
$arr[10,5] = 1.506398507577657
$arr[10,6] = 1.724921979692382
$arr[10,7] = 0.213150811894316
$arr[10,8] = 3.626748579528908
$arr[10,9] = 0.2184875378460163
$arr[11,0] = 0.171251766691924
$arr[11,1] = 0.1155300053563053
$arr[11,2] = 0.1472938330477156
$arr[11,3] = 2.426047418858531
$arr[11,4] = 1.715703619572707
$arr[11,5] = 1.507103509174527
$arr[11,6] = 1.725394096523821
$arr[11,7] = 0.2131305494075093
$arr[11,8] = 3.613416295767877
$arr[11,9] = 0.2184267032960818
$arr[12,0] = 0.1704162117927837
$arr[12,1] = 0.114808507135848
$arr[12,2] = 0.1469744114975207
$arr[12,3] = 2.428026599994112
$arr[12,4] = 1.71776359472139
$arr[12,5] = 1.509426003816557
$arr[12,6] = 1.726964113051132
$arr[12,7] = 0.2130668229534578
$arr[12,8] = 3.569943970788529
$arr[12,9] = 0.2182302648754941
$arr[13,0] = 0.1699050842887146
$arr[13,1] = 0.1143667189413264
$arr[13,2] = 0.1467796277555173
$arr[13,3] = 2.429270978806599
$arr[13,4] = 1.719052966058783
$arr[13,5] = 1.510866521389019
$arr[13,6] = 1.727949097860808
$arr[13,7] = 0.2130295975209719
$arr[13,8] = 3.54332144813975
$arr[13,9] = 0.2181114443582644
$arr[14,0] = 0.1669863966453136
$arr[14,1] = 0.1118373547164708
$arr[14,2] = 0.1456767756356179
$arr[14,3] = 2.436893061831825
$arr[14,4] = 1.726863796675019
$arr[14,5] = 1.51939601953724
$arr[14,6] = 1.733950420985309
$arr[14,7] = 0.2128440013173716
$arr[14,8] = 3.390844873720937
$arr[14,9] = 0.2174535983573804
$arr[15,0] = 0.1652051876710203
$arr[15,1] = 0.1102878101140163
$arr[15,2] = 0.1450121701229854
$arr[15,3] = 2.442006546217215
$arr[15,4] = 1.732031236527021
$arr[15,5] = 1.524873259864336
$arr[15,6] = 1.737949938960789
$arr[15,7] = 0.2127548891552209
$arr[15,8] = 3.297384183878307
$arr[15,9] = 0.2170706343720639
$arr[16,0] = 0.1641841257869316
$arr[16,1] = 0.1093973571579454
$arr[16,2] = 0.1446343109142312
$arr[16,3] = 2.445108263914918
$arr[16,4] = 1.735141189161382
$arr[16,5] = 1.528113440434893
$arr[16,6] = 1.740367003352858
$arr[16,7] = 0.2127127292454105
$arr[16,8] = 3.243657690921111
$arr[16,9] = 0.2168579435093037
$arr[17,0] = 0.1638390049745624
$arr[17,1] = 0.1090960049394596
$arr[17,2] = 0.1445071316365656
$arr[17,3] = 2.446186009288738
$arr[17,4] = 1.736217799939411
$arr[17,5] = 1.5292259266988
$arr[17,6] = 1.741205397073536
$arr[17,7] = 0.2127000168359743
$arr[17,8] = 3.225471976975541
$arr[17,9] = 0.2167872326812912
$arr[18,0] = 0.1653944445948952
$arr[18,1] = 0.1104526789726208
$arr[18,2] = 0.145082462902419
$arr[18,3] = 2.441445582772289
$arr[18,4] = 1.731466887330967
$arr[18,5] = 1.524280900396363
$arr[18,6] = 1.737512106854652
$arr[18,7] = 0.2127634340253834
$arr[18,8] = 3.307330183417662
$arr[18,9] = 0.2171106172401736
$arr[19,0] = 0.1706614941724638
$arr[19,1] = 0.1150203974301434
$arr[19,2] = 0.147068051848759
$arr[19,3] = 2.427438581182514
$arr[19,4] = 1.717152781328664
$arr[19,5] = 1.508740095248143
$arr[19,6] = 1.726498103527575
$arr[19,7] = 0.2130851639990894
$arr[19,8] = 3.582711697331831
$arr[19,9] = 0.2182876502792794
$arr[20,0] = 0.1741340369829345
$arr[20,1] = 0.118012441299129
$arr[20,2] = 0.148404766645438
$arr[20,3] = 2.419719087834551
$arr[20,4] = 1.709031140218457
$arr[20,5] = 1.499385186094159
$arr[20,6] = 1.720342378578884
$arr[20,7] = 0.2133763968709985
$arr[20,8] = 3.76293884676187
$arr[20,9] = 0.2191242330568741
$arr[21,0] = 0.1722779309697415
$arr[21,1] = 0.1164149295667585
$arr[21,2] = 0.1476877763607014
$arr[21,3] = 2.42370759235682
$arr[21,4] = 1.713252633602707
$arr[21,5] = 1.504304600741762
$arr[21,6] = 1.723532224199992
$arr[21,7] = 0.2132135534101636
$arr[21,8] = 3.666726089362953
$arr[21,9] = 0.2186715812952542
$arr[22,0] = 0.1653088722573557
$arr[22,1] = 0.1103781405093684
$arr[22,2] = 0.1450506703550083
$arr[22,3] = 2.441698690074844
$arr[22,4] = 1.731721596322274
$arr[22,5] = 1.524548421921978
$arr[22,6] = 1.737709684204681
$arr[22,7] = 0.2127595426320781
$arr[22,8] = 3.302833580351489
$arr[22,9] = 0.2170925176844065
$arr[23,0] = 0.157923462767485
$arr[23,1] = 0.1038974446565675
$arr[23,2] = 0.1423746213999664
$arr[23,3] = 2.467241206042118
$arr[23,4] = 1.756906992948956
$arr[23,5] = 1.549817328934239
$arr[23,6] = 1.757459756427679
$arr[23,7] = 0.2126175521153044
$arr[23,8] = 2.911456778013076
$arr[23,9] = 0.2156791332890791

$ws.Range("C2:L25").Value = $arr
